$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values --------------------------------------------------------------
# B1 and A2 hold numeric 0, B2 holds the label (goes into sharedStrings.xml)
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Formatting ------------------------------------------------------------
# B1 and A2 get: bold font, centered/top aligned, thin box border.
# Build the full style on B1 first (keeps styles.xml clean: one font, one
# border, one cellXf), then clone it onto A2 via copy/paste-special so the
# second range reuses the already-resolved style instead of generating a
# throwaway intermediate cellXf.
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight = 2            # xlThin

$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
